$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value (B11) and total corrected marks (B12)
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 95

# Update the "correct/total" summary text (E12)
$ws.Range("E12").Value = "95/140"
